$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for the 2022-Q3 quarter
#    right under the header, pushing all the other quarters down by
#    one row (2022-Q2 -> row3, 2022-Q1 -> row4, ... 2020-Q4 -> row9).
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# The auto "format copy" that Insert() performs leaves stray styling
# on the data cells of the new row - clear it so it matches the
# plain (unstyled) look of the other data rows.
$summary.Range("B2:D2").ClearFormats()

# Column A carries the bold/bordered "index" style - copy it from the
# row right below (which just got shifted down from the old row 2).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.79

# ------------------------------------------------------------------
# 2) Brand-new "2022-Q3" fund-detail sheet, inserted right after
#    "总计" (i.e. before the sheet that used to be first, "2022-Q2").
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # headers start at column B
    $cell = $q3.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

$idxCell = $q3.Cells.Item(2, 1)
$idxCell.Value = 0
$idxCell.Font.Bold = $true
$idxCell.HorizontalAlignment = -4108
$idxCell.VerticalAlignment = -4160
$idxCell.Borders.Item(7).LineStyle = 1
$idxCell.Borders.Item(8).LineStyle = 1
$idxCell.Borders.Item(9).LineStyle = 1
$idxCell.Borders.Item(10).LineStyle = 1

$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "010714"
$q3.Range("C2").Value = "东方红远见价值混合"

$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "15.24"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "94.15"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "5.17"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.7879"

$q3.Range("H2").Value = 4

# Keep "2020-Q4" (the last sheet) as the selected/active tab, matching
# the pre-edit state (adding a sheet makes it active by default).
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()
